$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Cells.Item(42, 8).Value = 474.66666
$ws.Cells.Item(42, 9).Value = 34.4
$ws.Cells.Item(42, 10).Value = 1025
$ws.Cells.Item(42, 11).Value = 103.2
$ws.Cells.Item(42, 12).Value = 3075
$ws.Cells.Item(42, 13).Value = 126.8
$ws.Cells.Item(42, 14).Value = -3535
# Row 55
$ws.Cells.Item(55, 8).Value = 84.2
$ws.Cells.Item(55, 9).Value = 84.2
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 84.2
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = 129.8
$ws.Cells.Item(55, 14).ClearContents()
# Row 129
$ws.Cells.Item(129, 8).Value = 1059.5172
$ws.Cells.Item(129, 10).Value = 1179.96
$ws.Cells.Item(129, 12).Value = 3539.88
$ws.Cells.Item(129, 14).Value = -13539.88

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5616.5713
$ws.Cells.Item(32, 9).Value = 3226.7166
$ws.Cells.Item(32, 11).Value = 3226.7166
$ws.Cells.Item(32, 13).Value = -2939.7166
# Row 61
$ws.Cells.Item(61, 8).Value = 1977.091
$ws.Cells.Item(61, 9).Value = 1764.3334
$ws.Cells.Item(61, 10).Value = 2544.4443
$ws.Cells.Item(61, 11).Value = 1764.3334
$ws.Cells.Item(61, 12).Value = 2544.4443
$ws.Cells.Item(61, 13).Value = -1552.3334
$ws.Cells.Item(61, 14).Value = -2968.4443
# Row 74
$ws.Cells.Item(74, 8).Value = 38224.816
$ws.Cells.Item(74, 9).Value = 51073
$ws.Cells.Item(74, 10).Value = 1515.7142
$ws.Cells.Item(74, 11).Value = 51073
$ws.Cells.Item(74, 12).Value = 1515.7142
$ws.Cells.Item(74, 13).Value = -50199
$ws.Cells.Item(74, 14).Value = -3263.7142
# Row 77
$ws.Cells.Item(77, 8).Value = 38224.816
$ws.Cells.Item(77, 9).Value = 51073
$ws.Cells.Item(77, 10).Value = 1515.7142
$ws.Cells.Item(77, 11).Value = 255365
$ws.Cells.Item(77, 12).Value = 7578.571
$ws.Cells.Item(77, 13).Value = -250997
$ws.Cells.Item(77, 14).Value = -16314.571
# Row 110
$ws.Cells.Item(110, 8).Value = 1764.9166
$ws.Cells.Item(110, 9).Value = 1863.5454
$ws.Cells.Item(110, 10).Value = 680
$ws.Cells.Item(110, 11).Value = 1863.5454
$ws.Cells.Item(110, 12).Value = 680
$ws.Cells.Item(110, 13).Value = 181.4546
$ws.Cells.Item(110, 14).Value = -4770
# Row 122
$ws.Cells.Item(122, 8).Value = 1800.25
$ws.Cells.Item(122, 9).Value = 1216.1538
$ws.Cells.Item(122, 10).Value = 2490.5454
$ws.Cells.Item(122, 11).Value = 3648.4614
$ws.Cells.Item(122, 12).Value = 7471.6362
$ws.Cells.Item(122, 13).Value = -1198.4614
$ws.Cells.Item(122, 14).Value = -12371.6362
# Row 136
$ws.Cells.Item(136, 8).Value = 1977.091
$ws.Cells.Item(136, 9).Value = 1764.3334
$ws.Cells.Item(136, 10).Value = 2544.4443
$ws.Cells.Item(136, 11).Value = 5293.0002
$ws.Cells.Item(136, 12).Value = 7633.3329
$ws.Cells.Item(136, 13).Value = -2743.0002
$ws.Cells.Item(136, 14).Value = -12733.3329

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 926.8857400000001
$ws.Cells.Item(20, 9).Value = 606.8
$ws.Cells.Item(20, 11).Value = 606.8
$ws.Cells.Item(20, 13).Value = -359.8
# Row 99
$ws.Cells.Item(99, 8).Value = 1242.7142
$ws.Cells.Item(99, 9).Value = 1274.75
$ws.Cells.Item(99, 10).Value = 1200
$ws.Cells.Item(99, 11).Value = 1274.75
$ws.Cells.Item(99, 12).Value = 1200
$ws.Cells.Item(99, 13).Value = 223.25
$ws.Cells.Item(99, 14).Value = -4196
# Row 112
$ws.Cells.Item(112, 8).Value = 28333.334
$ws.Cells.Item(112, 10).Value = 28333.334
$ws.Cells.Item(112, 12).Value = 28333.334
$ws.Cells.Item(112, 14).Value = -31287.334
# Row 134
$ws.Cells.Item(134, 8).Value = 5519.636
$ws.Cells.Item(134, 9).Value = 7103.1665
$ws.Cells.Item(134, 10).Value = 4925.8125
$ws.Cells.Item(134, 11).Value = 21309.4995
$ws.Cells.Item(134, 12).Value = 14777.4375
$ws.Cells.Item(134, 13).Value = -18774.4995
$ws.Cells.Item(134, 14).Value = -19847.4375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 27028842
$ws.Cells.Item(31, 9).Value = 52632510
$ws.Cells.Item(31, 10).Value = 2746.6667
$ws.Cells.Item(31, 11).Value = 52632510
$ws.Cells.Item(31, 12).Value = 2746.6667
$ws.Cells.Item(31, 13).Value = -52632215
$ws.Cells.Item(31, 14).Value = -3336.6667
# Row 34
$ws.Cells.Item(34, 8).Value = 27028842
$ws.Cells.Item(34, 9).Value = 52632510
$ws.Cells.Item(34, 10).Value = 2746.6667
$ws.Cells.Item(34, 11).Value = 52632510
$ws.Cells.Item(34, 12).Value = 2746.6667
$ws.Cells.Item(34, 13).Value = -52632308
$ws.Cells.Item(34, 14).Value = -3150.6667
# Row 132
$ws.Cells.Item(132, 8).Value = 2464.641
$ws.Cells.Item(132, 9).Value = 1889.2593
$ws.Cells.Item(132, 11).Value = 5667.7779
$ws.Cells.Item(132, 13).Value = -3137.7779
# Row 134
$ws.Cells.Item(134, 8).Value = 27501526
$ws.Cells.Item(134, 9).Value = 3449886.2
$ws.Cells.Item(134, 11).Value = 10349658.6
$ws.Cells.Item(134, 13).Value = -10347123.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Cells.Item(3, 8).Value = 4546.6665
$ws.Cells.Item(3, 10).Value = 6500
$ws.Cells.Item(3, 12).Value = 19500
$ws.Cells.Item(3, 14).Value = -19724
# Row 113
$ws.Cells.Item(113, 8).Value = 1515662.2
$ws.Cells.Item(113, 9).Value = 3788380.8
$ws.Cells.Item(113, 10).Value = 516.6667
$ws.Cells.Item(113, 11).Value = 11365142.4
$ws.Cells.Item(113, 12).Value = 1550.0001
$ws.Cells.Item(113, 13).Value = -11362972.4
$ws.Cells.Item(113, 14).Value = -5890.0001
# Row 122
$ws.Cells.Item(122, 8).Value = 975.3214
$ws.Cells.Item(122, 9).Value = 518.4706
$ws.Cells.Item(122, 10).Value = 1681.3636
$ws.Cells.Item(122, 11).Value = 4666.2354
$ws.Cells.Item(122, 12).Value = 15132.2724
$ws.Cells.Item(122, 13).Value = -2216.2354
$ws.Cells.Item(122, 14).Value = -20032.2724
# Row 123
$ws.Cells.Item(123, 8).Value = 8933.333000000001
$ws.Cells.Item(123, 9).Value = 12000
$ws.Cells.Item(123, 10).Value = 7400
$ws.Cells.Item(123, 11).Value = 36000
$ws.Cells.Item(123, 12).Value = 22200
$ws.Cells.Item(123, 13).Value = -33550
$ws.Cells.Item(123, 14).Value = -27100
# Row 131
$ws.Cells.Item(131, 8).Value = 852.77
$ws.Cells.Item(131, 9).Value = 588.1667
$ws.Cells.Item(131, 10).Value = 888.8523
$ws.Cells.Item(131, 11).Value = 1764.5001
$ws.Cells.Item(131, 12).Value = 2666.5569
$ws.Cells.Item(131, 13).Value = 3275.4999
$ws.Cells.Item(131, 14).Value = -12746.5569

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Cells.Item(4, 8).Value = 6998
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 6998
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 6998
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).Value = -7222
# Row 97
$ws.Cells.Item(97, 8).Value = 1445.2941
$ws.Cells.Item(97, 9).Value = 1046.25
$ws.Cells.Item(97, 10).Value = 1800
$ws.Cells.Item(97, 11).Value = 1046.25
$ws.Cells.Item(97, 12).Value = 1800
$ws.Cells.Item(97, 13).Value = -550.25
$ws.Cells.Item(97, 14).Value = -2792
# Row 113
$ws.Cells.Item(113, 8).Value = 1557.1428
$ws.Cells.Item(113, 9).Value = 1180
$ws.Cells.Item(113, 10).Value = 2500
$ws.Cells.Item(113, 11).Value = 1180
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 13).Value = 990
$ws.Cells.Item(113, 14).Value = -6840
# Row 122
$ws.Cells.Item(122, 8).Value = 55180.156
$ws.Cells.Item(122, 9).Value = 68634.87
$ws.Cells.Item(122, 10).Value = 4725
$ws.Cells.Item(122, 11).Value = 205904.61
$ws.Cells.Item(122, 12).Value = 14175
$ws.Cells.Item(122, 13).Value = -203454.61
$ws.Cells.Item(122, 14).Value = -19075
# Row 132
$ws.Cells.Item(132, 8).Value = 2981.6943
$ws.Cells.Item(132, 9).Value = 2891.6333
$ws.Cells.Item(132, 11).Value = 8674.8999
$ws.Cells.Item(132, 13).Value = -6144.8999
# Row 135
$ws.Cells.Item(135, 8).Value = 33000
$ws.Cells.Item(135, 10).Value = 33000
$ws.Cells.Item(135, 12).Value = 33000
$ws.Cells.Item(135, 14).Value = -43140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4406.2
$ws.Cells.Item(7, 9).Value = 4407.125
$ws.Cells.Item(7, 10).Value = 4402.5
$ws.Cells.Item(7, 11).Value = 4407.125
$ws.Cells.Item(7, 12).Value = 4402.5
$ws.Cells.Item(7, 13).Value = -4295.125
$ws.Cells.Item(7, 14).Value = -4626.5
# Row 13
$ws.Cells.Item(13, 8).Value = 1000
$ws.Cells.Item(13, 9).Value = 1000
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 1000
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 14).ClearContents()
$ws.Cells.Item(13, 13).Value = -860
# Row 43
$ws.Cells.Item(43, 8).Value = 9006
$ws.Cells.Item(43, 9).Value = 12
$ws.Cells.Item(43, 11).Value = 12
$ws.Cells.Item(43, 13).Value = 181
# Row 61
$ws.Cells.Item(61, 8).Value = 1217.0731
$ws.Cells.Item(61, 9).Value = 1132.8334
$ws.Cells.Item(61, 10).Value = 1446.8182
$ws.Cells.Item(61, 11).Value = 1132.8334
$ws.Cells.Item(61, 12).Value = 1446.8182
$ws.Cells.Item(61, 13).Value = -930.8334
$ws.Cells.Item(61, 14).Value = -1850.8182
# Row 68
$ws.Cells.Item(68, 8).Value = 10180.154
$ws.Cells.Item(68, 9).Value = 50651
$ws.Cells.Item(68, 10).Value = 2821.818
$ws.Cells.Item(68, 11).Value = 50651
$ws.Cells.Item(68, 12).Value = 2821.818
$ws.Cells.Item(68, 13).Value = -49902
$ws.Cells.Item(68, 14).Value = -4319.818
# Row 71
$ws.Cells.Item(71, 8).Value = 10180.154
$ws.Cells.Item(71, 9).Value = 50651
$ws.Cells.Item(71, 10).Value = 2821.818
$ws.Cells.Item(71, 11).Value = 253255
$ws.Cells.Item(71, 12).Value = 14109.09
$ws.Cells.Item(71, 13).Value = -249511
$ws.Cells.Item(71, 14).Value = -21597.09
# Row 113
$ws.Cells.Item(113, 8).Value = 1217.0731
$ws.Cells.Item(113, 9).Value = 1132.8334
$ws.Cells.Item(113, 10).Value = 1446.8182
$ws.Cells.Item(113, 11).Value = 1132.8334
$ws.Cells.Item(113, 12).Value = 1446.8182
$ws.Cells.Item(113, 13).Value = 1037.1666
$ws.Cells.Item(113, 14).Value = -5786.8182
# Row 126
$ws.Cells.Item(126, 8).Value = 4406.2
$ws.Cells.Item(126, 9).Value = 4407.125
$ws.Cells.Item(126, 10).Value = 4402.5
$ws.Cells.Item(126, 11).Value = 13221.375
$ws.Cells.Item(126, 12).Value = 13207.5
$ws.Cells.Item(126, 13).Value = -10751.375
$ws.Cells.Item(126, 14).Value = -18147.5
# Row 136
$ws.Cells.Item(136, 8).Value = 8131352
$ws.Cells.Item(136, 9).Value = 11495571
$ws.Cells.Item(136, 10).Value = 1155
$ws.Cells.Item(136, 11).Value = 34486713
$ws.Cells.Item(136, 12).Value = 3465
$ws.Cells.Item(136, 13).Value = -34484163
$ws.Cells.Item(136, 14).Value = -8565

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 90909820
$ws.Cells.Item(126, 9).Value = 100000700
$ws.Cells.Item(126, 10).Value = 1000
$ws.Cells.Item(126, 11).Value = 300002100
$ws.Cells.Item(126, 12).Value = 3000
$ws.Cells.Item(126, 13).Value = -299999630
$ws.Cells.Item(126, 14).Value = -7940
# Row 132
$ws.Cells.Item(132, 8).Value = 7696842
$ws.Cells.Item(132, 9).Value = 10531100
$ws.Cells.Item(132, 10).Value = 3856.8572
$ws.Cells.Item(132, 11).Value = 31593300
$ws.Cells.Item(132, 12).Value = 11570.5716
$ws.Cells.Item(132, 13).Value = -31590770
$ws.Cells.Item(132, 14).Value = -16630.5716
